# Apply hybrid bold + color highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) within specific bullet
# paragraphs, matching the target OOXML diff.

function ToWordColor($hex) {
    # Word's Range.Font.Color is an OLE_COLOR (0x00BBGGRR); convert an
    # "RRGGBB" hex string into that representation.
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $b * 65536 + $g * 256 + $r
}

$d = $word.ActiveDocument
$highlightColor = ToWordColor "2C3E50"
$pm = [string][char]13
$bullet = [string][char]0x2022
$plusMinus = [string][char]0x00B1

# Find the paragraph whose text (sans trailing paragraph mark) exactly
# equals $ExactText and apply bold + highlight color to each substring
# listed in $Metrics (applied left-to-right, each search resuming after
# the previous match so repeated values - e.g. two "87%" across two
# different paragraphs - are handled correctly within a single call).
function Set-MetricHighlights {
    param(
        [string]$ExactText,
        [string[]]$Metrics
    )

    $count = $d.Paragraphs.Count
    $target = $null
    for ($i = 1; $i -le $count; $i++) {
        $para = $d.Paragraphs.Item($i)
        $t = $para.Range.Text.TrimEnd($pm)
        if ($t -eq $ExactText) {
            $target = $para
        }
    }

    if ($target -eq $null) {
        Write-Host "Paragraph not found for:" $ExactText
        return
    }

    $range = $target.Range
    $paraStart = $range.Start
    $full = $range.Text.TrimEnd($pm)

    $searchFrom = 0
    foreach ($metric in $Metrics) {
        $idx = $full.IndexOf($metric, $searchFrom)
        if ($idx -lt 0) {
            Write-Host "Metric not found:" $metric
            continue
        }
        $start = $paraStart + $idx
        $end = $start + $metric.Length
        $sub = $d.Range($start, $end)
        $sub.Font.Bold = 1
        $sub.Font.Color = $highlightColor
        $searchFrom = $idx + $metric.Length
    }
}

Set-MetricHighlights `
    ($bullet + " Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%") `
    @("23%", "64%")

Set-MetricHighlights `
    ($bullet + " Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from " + $plusMinus + "4.2% to " + $plusMinus + "2.1%") `
    @("87%", "71%", ($plusMinus + "4.2%"), ($plusMinus + "2.1%"))

Set-MetricHighlights `
    ($bullet + " Wrote RFP and analyzed bids from 1,200 vendors for research platform development") `
    @("1,200")

Set-MetricHighlights `
    ($bullet + " Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+") `
    @('$400M', '$1B')

Set-MetricHighlights `
    ($bullet + " Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M") `
    @("73.5%", '$4.7M')

Set-MetricHighlights `
    ($bullet + " Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%") `
    @("87%", "71%")
